$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.829.60"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "2.811.03"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("E4").Value = "  +0.16%  "
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "351.10"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.92%  "
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "113.19"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +4.66%  "
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.561"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +2.25%  "
$ws.Range("E8").Value = "  +0.04%  "
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.619"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +6.17%  "
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "40.33"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("E11").Value = "  -0.84%  "
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0844"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +1.22%  "
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "19.92"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -0.02%  "
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.80"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +3.66%  "
$ws.Range("D15").Value = "3.259.55"
$ws.Range("E15").Value = "  +1.87%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.968"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +4.30%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.810.46"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").Value = "51.924.82"
$ws.Range("E18").Value = "  +1.50%  "
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.38"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +9.63%  "
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.61"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.72%  "
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "13.56"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +3.21%  "
$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("E22").Value = "  +1.36%  "
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "70.61"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.35%  "
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "268.87"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +1.51%  "
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.76"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.86%  "
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "26.21"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.90%  "
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.52"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +3.42%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "38.71"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +11.03%  "
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.27"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +1.12%  "
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.16"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +1.34%  "
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0906"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +9.03%  "
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.69"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +2.65%  "
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0453"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +2.38%  "
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.02%  "
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "19.02"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +4.66%  "
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.21"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("E40").Value = "  +3.04%  "
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.59"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +2.91%  "
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.116"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +2.02%  "
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "22.42"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +1.71%  "
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.24"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +2.03%  "
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "120.83"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.26%  "
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.51"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +8.44%  "
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.49"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +10.05%  "
$ws.Range("D48").Value = "2.139.74"
$ws.Range("E48").Value = "  +2.53%  "
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.02"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +11.12%  "
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.226"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +19.16%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0322"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +15.78%  "
